{"js": "// Replace the date line and each three-digit-by-one-digit division\n// problem text with the updated values from the new day's worksheet.\nconst replacements = [[\"2024-07-15 Monday\", \"2024-07-16 Tuesday\"], [\"641\u00f73=\", \"238\u00f79=\"], [\"464\u00f74=\", \"358\u00f79=\"], [\"723\u00f76=\", \"888\u00f72=\"], [\"285\u00f76=\", \"898\u00f72=\"], [\"612\u00f75=\", \"888\u00f75=\"], [\"417\u00f76=\", \"389\u00f79=\"], [\"913\u00f79=\", \"389\u00f77=\"], [\"814\u00f73=\", \"938\u00f77=\"], [\"143\u00f77=\", \"169\u00f78=\"], [\"810\u00f76=\", \"119\u00f76=\"], [\"466\u00f74=\", \"277\u00f74=\"], [\"931\u00f72=\", \"157\u00f73=\"], [\"178\u00f78=\", \"649\u00f76=\"], [\"967\u00f74=\", \"256\u00f76=\"], [\"212\u00f79=\", \"269\u00f73=\"], [\"600\u00f72=\", \"376\u00f77=\"], [\"509\u00f72=\", \"510\u00f78=\"], [\"900\u00f72=\", \"100\u00f79=\"], [\"696\u00f77=\", \"250\u00f74=\"], [\"225\u00f73=\", \"335\u00f72=\"], [\"954\u00f76=\", \"496\u00f76=\"], [\"676\u00f75=\", \"140\u00f74=\"], [\"561\u00f72=\", \"881\u00f73=\"], [\"496\u00f74=\", \"479\u00f79=\"], [\"834\u00f72=\", \"103\u00f76=\"]];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# wdReplaceAll = 2, wdFindContinue = 1\n$replacements = @(\n    @(\"2024-07-15 Monday\", \"2024-07-16 Tuesday\"),\n    @(\"641\u00f73=\", \"238\u00f79=\"),\n    @(\"464\u00f74=\", \"358\u00f79=\"),\n    @(\"723\u00f76=\", \"888\u00f72=\"),\n    @(\"285\u00f76=\", \"898\u00f72=\"),\n    @(\"612\u00f75=\", \"888\u00f75=\"),\n    @(\"417\u00f76=\", \"389\u00f79=\"),\n    @(\"913\u00f79=\", \"389\u00f77=\"),\n    @(\"814\u00f73=\", \"938\u00f77=\"),\n    @(\"143\u00f77=\", \"169\u00f78=\"),\n    @(\"810\u00f76=\", \"119\u00f76=\"),\n    @(\"466\u00f74=\", \"277\u00f74=\"),\n    @(\"931\u00f72=\", \"157\u00f73=\"),\n    @(\"178\u00f78=\", \"649\u00f76=\"),\n    @(\"967\u00f74=\", \"256\u00f76=\"),\n    @(\"212\u00f79=\", \"269\u00f73=\"),\n    @(\"600\u00f72=\", \"376\u00f77=\"),\n    @(\"509\u00f72=\", \"510\u00f78=\"),\n    @(\"900\u00f72=\", \"100\u00f79=\"),\n    @(\"696\u00f77=\", \"250\u00f74=\"),\n    @(\"225\u00f73=\", \"335\u00f72=\"),\n    @(\"954\u00f76=\", \"496\u00f76=\"),\n    @(\"676\u00f75=\", \"140\u00f74=\"),\n    @(\"561\u00f72=\", \"881\u00f73=\"),\n    @(\"496\u00f74=\", \"479\u00f79=\"),\n    @(\"834\u00f72=\", \"103\u00f76=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $found = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        throw \"Text not found: $oldText\"\n    }\n}\n"}
